$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.29679999999999
$ws.Range("C7").Value = -13.33369999999999
$ws.Range("A10").Value = -21.67489999999999
$ws.Range("A12").Value = -21.56669999999999
$ws.Range("C15").Value = -14.29339999999998
$ws.Range("A18").Value = -22.29080000000001
$ws.Range("D18").Value = -8.173599999999993
$ws.Range("D19").Value = -8.890399999999993
$ws.Range("C20").Value = -11.9074
$ws.Range("D27").Value = -9.0063
$ws.Range("C29").Value = -11.43090000000001
$ws.Range("C30").Value = -12.79189999999999
$ws.Range("C31").Value = -13.3266
$ws.Range("A37").Value = -20.62340000000001
$ws.Range("C40").Value = -13.55050000000001
$ws.Range("D42").Value = -8.688899999999993
$ws.Range("D44").Value = -7.4522
$ws.Range("D47").Value = -7.554199999999997
$ws.Range("A55").Value = -22.2976
$ws.Range("D58").Value = -8.441899999999993
$ws.Range("A68").Value = -21.632
$ws.Range("C68").Value = -11.6218
$ws.Range("D73").Value = -7.5914
$ws.Range("C76").Value = -12.0884
$ws.Range("A77").Value = -20.7957
$ws.Range("A78").Value = -20.46299999999998
$ws.Range("C87").Value = -13.80019999999999
$ws.Range("C88").Value = -13.37649999999999
$ws.Range("D95").Value = -7.7357
$ws.Range("C96").Value = -12.9976
$ws.Range("C98").Value = -12.3704
$ws.Range("C101").Value = -13.74090000000001
$ws.Range("D101").Value = -8.070799999999997
$ws.Range("C102").Value = -13.30270000000001
